# Integración con perfiles definidos
# Updates the "Casino" permission test-data row (row 2) with a new,
# successful test run: new URL, account/user, date/time stamps, OS /
# browser / resolution info and a new "Exitoso" outcome.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to hold the value as literal text (matches the
    # original sheet, where every data cell in the row is a shared
    # string) instead of letting Excel auto-detect numbers/dates/times,
    # then drop back to General so no stray explicit style sticks to
    # the cell.
    $range.Clear()
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

# B2 previously carried a hyperlink (and its "visited link" style) to
# the old support URL - remove both before writing the new URL so the
# cell goes back to plain/default formatting.
$ws.Range("B2").Hyperlinks.Delete()

Set-TextValue $ws.Range("B2") "http://10.0.74.5/index/home"
Set-TextValue $ws.Range("C2") "4fku01"
Set-TextValue $ws.Range("D2") "123"
Set-TextValue $ws.Range("E2") "10/10/2019"
Set-TextValue $ws.Range("F2") "14:12:45.637"
Set-TextValue $ws.Range("G2") "10/10/2019"
Set-TextValue $ws.Range("H2") "14:13:09.318"
Set-TextValue $ws.Range("I2") "Windows Server 2016"
Set-TextValue $ws.Range("J2") "Firefox 69.0.2"
Set-TextValue $ws.Range("K2") "1920x1080"
Set-TextValue $ws.Range("L2") "Casino"
Set-TextValue $ws.Range("M2") "Exitoso"
Set-TextValue $ws.Range("N2") "Permiso Casino es correctamente accesible para el usuario"
